$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New token-table content (rows 2-15), replacing the previous "printf" example
# with a new "int main() { int 1stPlace = 100; return 0; }" example.
$data = @(
    @("int",      "keyword",    "<INT_TK>"),
    @("main",     "keyword",    "<MAIN_TK>"),
    @("(",        "symbol",     "<PHARANTESES1_TK>"),
    @(")",        "symbol",     "<PHARANTESES2_TK>"),
    @("{",        "symbol",     "<BRACKET1_TK>"),
    @("int",      "keyword",    "<INT_TK>"),
    @("1stPlace", "identifier", "<ID_TK, 1>"),
    @("=",        "operators",  "<ASIGN_TK>"),
    @("100",      "integer",    "<INT_CONST>"),
    @(";",        "symbol",     "<SEMICOLON_TK>"),
    @("return",   "keyword",    "<RETURN_TK>"),
    @("0",        "integer",    "<INT_CONST>"),
    @(";",        "symbol",     "<SEMICOLON_TK>"),
    @("}",        "symbol",     "<BRACKET2_TK>")
)

$rowIndex = 2
foreach ($entry in $data) {
    $codeValue = $entry[0]
    $cellA = $ws.Cells.Item($rowIndex, 1)
    if ($codeValue.StartsWith("=")) {
        # Prevent Excel from parsing the code sample (e.g. "=") as a formula,
        # then reset the style so no extra quote-prefix formatting lingers.
        $cellA.Value = "'" + $codeValue
        $cellA.Style = "Normal"
    }
    elseif ($codeValue -match '^[0-9]+$') {
        # Keep numeric-looking code samples (e.g. "100", "0") stored as text,
        # matching how every other cell in this column is stored.
        $cellA.NumberFormat = "@"
        $cellA.Value = $codeValue
        $cellA.Style = "Normal"
    }
    else {
        $cellA.Value = $codeValue
    }
    $ws.Cells.Item($rowIndex, 2).Value = $entry[1]
    $ws.Cells.Item($rowIndex, 3).Value = $entry[2]
    $rowIndex++
}

# Remove the old trailing rows (16-18) that no longer exist in the table.
$ws.Range("A16:C18").Delete()
